$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.857.98'
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = '1.635.99'
$ws.Range("E3").Value = '  -1.26%  '
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.72'
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5031'
$ws.Range("E6").Value = '  -2.27%  '
$ws.Range("E7").Value = '  -0.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2573'
$ws.Range("E8").Value = '  -0.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06422'
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("E10").Value = '  -1.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07699'
$ws.Range("E11").Value = '  -1.12%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.244'
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.629.60'
$ws.Range("E13").Value = '  -1.71%  '
$ws.Range("D14").Value = '1.859.51'
$ws.Range("E14").Value = '  -1.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5454'
$ws.Range("E15").Value = '  -1.47%  '
$ws.Range("D16").Value = '0.0₅7948'
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").Value = '25.860.30'
$ws.Range("E18").Value = '  -1.34%  '
$ws.Range("E19").Value = '  -0.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '203.10'
$ws.Range("E20").Value = '  -3.87%  '
$ws.Range("E21").Value = '  -1.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.954'
$ws.Range("E22").Value = '  -1.11%  '
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.931'
$ws.Range("E25").Value = '  +11.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.08'
$ws.Range("E26").Value = '  -2.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1148'
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.717'
$ws.Range("E29").Value = '  -3.73%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.243'
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05022'
$ws.Range("E31").Value = '  -2.39%  '
$ws.Range("E32").Value = '  -2.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.184'
$ws.Range("E33").Value = '  -1.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.536'
$ws.Range("E34").Value = '  -2.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.359'
$ws.Range("E35").Value = '  -0.55%  '
$ws.Range("D36").Value = '1.171.56'
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.8948'
$ws.Range("E37").Value = '  -3.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.606'
$ws.Range("E38").Value = '  -5.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5619'
$ws.Range("E39").Value = '  -1.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01560'
$ws.Range("E40").Value = '  -2.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.541'
$ws.Range("E42").Value = '  -1.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.670'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8070'
$ws.Range("E44").Value = '  -3.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '99.59'
$ws.Range("E45").Value = '  -0.73%  '
$ws.Range("D46").Value = '1.771.66'
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").Value = '0.0₈115'
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4513'
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.003'
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '54.80'
$ws.Range("E50").Value = '  -1.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05051'
$ws.Range("E51").Value = '  -0.10%  '
